# Update countries & provincias Spain
# Applies the 23-Mar-2020 05:46 -> 06:16 data refresh to the "Pais" sheet:
#   - Tailandia, India, Mexico and Sri Lanka each overtake their neighbours in
#     "Casos totales" and move to the top of their local group; the other
#     countries in that group keep their own figures but shift down one row.
#   - The "Datos actualizados..." timestamp cell is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 06:16"

# Each entry: row number, country name, then Casos totales / Nuevos casos /
# Casos activos / Recuperados / Casos criticos / Muertes hoy / Muertes
$updates = @(
    @(33, "Tailandia", 721, 122,  52, 668,  7, 0,  1),
    @(34, "Crucero",   712,   0, 567, 137, 15, 0,  8),
    @(35, "Polonia",   634,   0,  13, 614,  3, 0,  7),
    @(36, "Chile",     632,   0,   8, 623,  7, 0,  1),
    @(37, "Finlandia", 626,   0,  10, 615, 12, 0,  1),
    @(38, "Grecia",    624,   0,  19, 590, 18, 0, 15),

    @(45, "India",     425,  29,  24, 393,  0, 1,  8),
    @(46, "Eslovenia", 414,   0,   0, 412, 12, 0,  2),

    @(54, "Mexico",    316,  65,   4, 310,  1, 0,  2),
    @(55, "Panama",    313,   0,   1, 309,  7, 0,  3),
    @(56, "Sudafrica", 274,   0,   2, 272,  0, 0,  0),
    @(57, "Argentina", 266,   0,  27, 235,  0, 0,  4),
    @(58, "Croacia",   254,   0,   5, 248,  5, 0,  1),

    @(90, "Sri Lanka",  86,   4,   3,  83,  2, 0,  0),
    @(91, "Camboya",    84,   0,   2,  82,  0, 0,  0)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 1).Value = $u[1]
    $ws.Cells.Item($row, 2).Value = $u[2]
    $ws.Cells.Item($row, 3).Value = $u[3]
    $ws.Cells.Item($row, 4).Value = $u[4]
    $ws.Cells.Item($row, 5).Value = $u[5]
    $ws.Cells.Item($row, 6).Value = $u[6]
    $ws.Cells.Item($row, 7).Value = $u[7]
    $ws.Cells.Item($row, 8).Value = $u[8]
}
